# Add a new "PatientLogin" worksheet with patient-role login test data,
# positioned after "patientInfo1" (so the tab order becomes
# UserLogin, patientInfo1, PatientLogin) and make it the active sheet.

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("UserLogin")
$wsPatientInfo = $wb.Worksheets.Item("patientInfo1")

$wsNew = $wb.Worksheets.Add($null, $wsPatientInfo)
$wsNew.Name = "PatientLogin"

# Write the data values first, then the header row, so that the shared
# string table is populated in the same order as the source data entry.
$wsNew.Range("A2").Value = "test"
$wsNew.Range("B2").Value = "rams@gmail.com"
$wsNew.Range("A1").Value = "patientpassword"
$wsNew.Range("B1").Value = "patinetEmail"

# The password-style value in B2 is actually an email -> add a mailto
# hyperlink on it, matching the pattern used on the UserLogin sheet.
$wsNew.Hyperlinks.Add($wsNew.Range("B2"), "mailto:rams@gmail.com")

# Match formatting used on the UserLogin sheet: header row style and the
# hyperlink cell style (copy formats only, so values/hyperlink stay put).
$wsLogin.Range("A1:B1").Copy()
$wsNew.Range("A1:B1").PasteSpecial(-4122)

$wsLogin.Range("B3").Copy()
$wsNew.Range("B2").PasteSpecial(-4122)

# Match the column widths used for the new sheet.
$wsNew.Columns.Item(1).ColumnWidth = 13.59
$wsNew.Columns.Item(2).ColumnWidth = 21.92

# Select B1 and make PatientLogin the active tab.
$wsNew.Range("B1").Select() | Out-Null
